# Re-stamp the seeded export rows with a fresh batch of generated ids /
# timestamps (as if the workbook had just been re-exported). Also fixes the
# "tid" header (should be "_tid", matching the other internal-field naming
# convention) on the Tags sheet so it can be parsed back in from Excel too.

$wb = $excel.ActiveWorkbook

# --- Defs ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Defs")
$ws.Range("A2").Value = "lgsauyu8-4goi"
$ws.Range("B2").Value = "2023-04-22T18:14:17.312Z"
$ws.Range("C2").Value = "lgsauyu9"

$ws.Range("A3").Value = "lgsauyu9-tbmb"
$ws.Range("B3").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C3").Value = "lgsauyu9"

$ws.Range("A4").Value = "lgsauyu9-m24n"
$ws.Range("B4").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C4").Value = "lgsauyu9"

# --- Point Defs -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Point Defs")
$ws.Range("A2").Value = "lgsauyu9-09r5"
$ws.Range("B2").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C2").Value = "lgsauyu9"

$ws.Range("A3").Value = "lgsauyu9-s0u5"
$ws.Range("B3").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C3").Value = "lgsauyu9"
$ws.Range("F3").Value = "0ksk"

$ws.Range("A4").Value = "lgsauyu9-v7tg"
$ws.Range("B4").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C4").Value = "lgsauyu9"

$ws.Range("A5").Value = "lgsauyu9-og4k"
$ws.Range("B5").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C5").Value = "lgsauyu9"

# --- Entry Base -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Base")
$ws.Range("A2").Value = "lgsauyu9-p2y4"
$ws.Range("B2").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C2").Value = "lgsauyu9"

$ws.Range("A3").Value = "lgsauyu9-3yeb"
$ws.Range("B3").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C3").Value = "lgsauyu9"
$ws.Range("F3").Value = "lgsauyul-0g07"
$ws.Range("G3").Value = "2023-04-22T13:14:17"

# --- Entry Points -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Points")
$ws.Range("A2").Value = "lgsauyu9-ydzh"
$ws.Range("B2").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C2").Value = "lgsauyu9"

$ws.Range("A3").Value = "lgsauyu9-mokn"
$ws.Range("B3").Value = "2023-04-22T18:14:17.313Z"
$ws.Range("C3").Value = "lgsauyu9"

# --- Tag Defs -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Tag Defs")
$ws.Range("A2").Value = "lgsauyul-xsvg"
$ws.Range("B2").Value = "2023-04-22T18:14:17.325Z"
$ws.Range("C2").Value = "lgsauyul"
$ws.Range("E2").Value = "05c2"

$ws.Range("A3").Value = "lgsauyul-h3kr"
$ws.Range("B3").Value = "2023-04-22T18:14:17.325Z"
$ws.Range("C3").Value = "lgsauyul"

$ws.Range("A4").Value = "lgsauyul-g8mm"
$ws.Range("B4").Value = "2023-04-22T18:14:17.325Z"
$ws.Range("C4").Value = "lgsauyul"

# --- Tags -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Tags")
$ws.Range("G1").Value = "_tid"

$ws.Range("A2").Value = "lgsauyul-gzp5"
$ws.Range("B2").Value = "2023-04-22T18:14:17.325Z"
$ws.Range("C2").Value = "lgsauyul"

$ws.Range("A3").Value = "lgsauyul-0jnx"
$ws.Range("B3").Value = "2023-04-22T18:14:17.325Z"
$ws.Range("C3").Value = "lgsauyul"
